# New weekly "Albahaca" price report row for Terminal La Palmera de La Serena.
# The source feed inserted one additional record between the existing rows
# 109 and 110 (by date order), which pushes every following record down by
# one row (old row 110 -> new row 111, ..., old row 168 -> new row 169).
#
# Insert a new physical row at 110 (EntireRow insert shifts 110..168 down
# to 111..169, matching the diff's dimension change A1:R168 -> A1:R169),
# then populate it with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(110).Insert()

$ws.Cells.Item(110, 1).Value  = 8
$ws.Cells.Item(110, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(110, 3).Value  = "Coquimbo"
$ws.Cells.Item(110, 4).Value  = 45001
$ws.Cells.Item(110, 5).Value  = 4
$ws.Cells.Item(110, 6).Value  = 100112052
$ws.Cells.Item(110, 7).Value  = "Albahaca"
$ws.Cells.Item(110, 8).Value  = "Sin especificar"
$ws.Cells.Item(110, 9).Value  = "Primera"
$ws.Cells.Item(110, 10).Value = 1200
$ws.Cells.Item(110, 11).Value = 2500
$ws.Cells.Item(110, 12).Value = 3000
$ws.Cells.Item(110, 13).Value = 2750
$ws.Cells.Item(110, 14).Value = "`$/docena de matas"
$ws.Cells.Item(110, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(110, 16).Value = 458
$ws.Cells.Item(110, 17).Value = 6
$ws.Cells.Item(110, 18).Value = "Hortaliza"
